# Automatic edit in sheet "Card13":
# Column A ("card") values for rows 3-7 and 9-13 are corrected from "2" to "13".
# (Row 2 and row 8 already contain "13" and are left untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card13")

# Use a cell that already holds the text value "13" (e.g. A2) as the copy
# source so the destination cells keep the same text (string) storage type
# as the rest of the column instead of being re-interpreted as numbers.
$source = $ws.Range("A2")

$rowsToFix = @(3, 4, 5, 6, 7, 9, 10, 11, 12, 13)
foreach ($r in $rowsToFix) {
    $target = $ws.Cells.Item($r, 1)
    $source.Copy($target)
}
